$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.207.59'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '1.839.48'
$ws.Range("E3").Value = '  -1.56%  '
$ws.Range("D4").Value = '''0.9999'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''232.60'
$ws.Range("E5").Value = '  -1.34%  '
$ws.Range("D6").Value = '''0.9999'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  -3.23%  '
$ws.Range("D8").Value = '''0.2706'
$ws.Range("E8").Value = '  -3.41%  '
$ws.Range("D9").Value = '''0.06271'
$ws.Range("E9").Value = '  -3.67%  '
$ws.Range("D10").Value = '1.834.03'
$ws.Range("E10").Value = '  -1.76%  '
$ws.Range("D11").Value = '''0.07415'
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '''16.09'
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("D13").Value = '''4.929'
$ws.Range("E13").Value = '  -2.93%  '
$ws.Range("D14").Value = '''83.70'
$ws.Range("E14").Value = '  -4.11%  '
$ws.Range("D15").Value = '''0.6205'
$ws.Range("E15").Value = '  -3.38%  '
$ws.Range("D16").Value = '30.135.37'
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").Value = '''0.9997'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Value = '''227.02'
$ws.Range("E18").Value = '  -1.95%  '
$ws.Range("D19").Value = '''0.000007295'
$ws.Range("E19").Value = '  -2.85%  '
$ws.Range("D20").Value = '''12.35'
$ws.Range("E20").Value = '  -4.94%  '
$ws.Range("D21").Value = '''0.9994'
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.075.49'
$ws.Range("E22").Value = '  -1.36%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '''4.879'
$ws.Range("E23").Value = '  -5.21%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '''5.844'
$ws.Range("E24").Value = '  -4.33%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '''9.200'
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '''164.10'
$ws.Range("E26").Value = '  -2.40%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''17.80'
$ws.Range("E27").Value = '  -3.42%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '''1.877'
$ws.Range("E28").Value = '  -1.97%  '
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").Value = '''0.1044'
$ws.Range("E29").Value = '  +1.19%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '''1.370'
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '''4.073'
$ws.Range("E31").Value = '  -4.64%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''3.790'
$ws.Range("E32").Value = '  -5.20%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '''0.04813'
$ws.Range("E33").Value = '  -3.43%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '''1.138'
$ws.Range("E34").Value = '  -3.35%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '''0.7094'
$ws.Range("E35").Value = '  -4.42%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '''2.692'
$ws.Range("E36").Value = '  -0.83%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.01885'
$ws.Range("E37").Value = '  -3.66%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '''2.648'
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '''0.8914'
$ws.Range("E39").Value = '  -3.10%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '''1.918'
$ws.Range("E40").Value = '  -6.55%  '
$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").Value = '''104.18'
$ws.Range("E41").Value = '  -1.63%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '''1.001'
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.535'
$ws.Range("E43").Value = '  -0.55%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '''0.4001'
$ws.Range("E44").Value = '  -4.84%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '''7.028'
$ws.Range("E45").Value = '  -2.58%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '''0.1194'
$ws.Range("E46").Value = '  -2.81%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '''59.74'
$ws.Range("E47").Value = '  -2.95%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''8.503'
$ws.Range("E48").Value = '  -4.46%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''32.74'
$ws.Range("E49").Value = '  -2.36%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.05508'
$ws.Range("E50").Value = '  -2.50%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '''1.355'
$ws.Range("E51").Value = '  -5.48%  '
